$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$apos = [char]0x2019

# ---------------------------------------------------------------------------
# 1) Paragraph "3. / A. One possible solution ..." -- fix "cage," -> "cage;"
#    and "then" -> "and then", dropping the w:proofErr wrappers around them.
# ---------------------------------------------------------------------------
$para3A = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*One possible solution*") {
        $para3A = $cand
        break
    }
}

$xmlPara3A = "<w:p $ns>" + `
    "<w:r><w:tab/><w:t xml:space=`"preserve`">A. One possible solution would be to put the bird in a </w:t></w:r>" + `
    "<w:r><w:t>cage;</w:t></w:r>" + `
    "<w:r><w:t xml:space=`"preserve`"> this will prevent the bird from eating the seed, or the cat from eating the bird.  The other solution might be to put the cat in a cage, </w:t></w:r>" + `
    "<w:r><w:t>and then</w:t></w:r>" + `
    "<w:r><w:t xml:space=`"preserve`"> the cat couldn${apos}t eat the bird and could be left alone with the seed. </w:t></w:r>" + `
    "</w:p>"
$para3A.Range.InsertXML($xmlPara3A)

# ---------------------------------------------------------------------------
# 2) Locate paragraph "4. " and insert six new answer paragraphs after it
#    (steps 4 and 5 of the problem-solving write-up), then relocate the
#    _GoBack bookmark to the very end of the new content.
# ---------------------------------------------------------------------------
$para4 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "4. ") {
        $para4 = $cand
        break
    }
}
$para4Index = $para4.Range.Start

# Create six fresh blank paragraphs right after paragraph "4. ".
$anchor = $d.Paragraphs.Item($para4.Index)
$anchor.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($para4.Index + 1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($para4.Index + 2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($para4.Index + 3)
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item($para4.Index + 4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item($para4.Index + 5)
$p5.Range.InsertParagraphAfter()

# Strip the _GoBack bookmark off paragraph "4. " (clean re-write of its text).
$xmlPara4 = "<w:p $ns><w:r><w:t xml:space=`"preserve`">4. </w:t></w:r></w:p>"
$d.Paragraphs.Item($para4.Index).Range.InsertXML($xmlPara4)

$xmlA = "<w:p $ns><w:r><w:tab/><w:t xml:space=`"preserve`">A.  I feel that each solution would successfully meet the goals.  </w:t></w:r></w:p>"
$d.Paragraphs.Item($para4.Index + 1).Range.InsertXML($xmlA)

$xmlB = "<w:p $ns><w:r><w:tab/><w:t xml:space=`"preserve`">B. The only issue with my solutions is the question of if the man owns a cage, if not, how is he transporting the parrot in the first place? Is he carrying the parrot on his shoulder? If so, then he could have the parrot on his shoulder, and the birdseed in the boat, and come back for the cat separate right? </w:t></w:r></w:p>"
$d.Paragraphs.Item($para4.Index + 2).Range.InsertXML($xmlB)

# Third new paragraph (blank spacer) is left empty on purpose.
$xmlEmpty = "<w:p $ns/>"
$d.Paragraphs.Item($para4.Index + 3).Range.InsertXML($xmlEmpty)

$xml5 = "<w:p $ns><w:r><w:t xml:space=`"preserve`">5. </w:t></w:r></w:p>"
$d.Paragraphs.Item($para4.Index + 4).Range.InsertXML($xml5)

$xml5A = "<w:p $ns><w:r><w:tab/><w:t xml:space=`"preserve`">A. My solution would be to put the parrot in a cage. This way the man could transport the parrot on it${apos}s own, leave it, come back and grab the seed and leave it with the parrot. The parrot being in the cage wouldn${apos}t be able to get to the seed, and the man would be able to head back safely grabbing the cat for transport. </w:t></w:r></w:p>"
$d.Paragraphs.Item($para4.Index + 5).Range.InsertXML($xml5A)

# Final paragraph carries the relocated _GoBack bookmark at its tail end.
$xml5B = "<w:p $ns><w:r><w:tab/><w:t xml:space=`"preserve`">B.  I am not sure a drawing is necessary for this problem. </w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$d.Paragraphs.Item($para4.Index + 6).Range.InsertXML($xml5B)

Write-Output "done"
